$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 2
$ws.Range("E2").Value = 1
$ws.Range("F2").Value = 0.3333333333333333
$ws.Range("G2").Value = 0.05881766666666666
$ws.Range("H2").Value = 0.176453
$ws.Range("M2").Value = 35.991783
$ws.Range("N2").Value = 107.975349
$ws.Range("O2").Value = 0.3909505149237033
$ws.Range("P2").Value = 0.3909505149237033
$ws.Range("Q2").Value = 2.116952695233
$ws.Range("R2").Value = 19.052574257097
$ws.Range("S2").Value = 0.3909505149237033
$ws.Range("T2").Value = 0.3909505149237033

# Row 3
$ws.Range("E3").Value = 1
$ws.Range("F3").Value = 0.3333333333333333
$ws.Range("G3").Value = 0.05881766666666666
$ws.Range("H3").Value = 0.176453
$ws.Range("O3").Value = 0.537552751174421
$ws.Range("P3").Value = 0.537552751174421
$ws.Range("Q3").Value = 2.910787176353222
$ws.Range("R3").Value = 26.197084587179
$ws.Range("S3").Value = 0.537552751174421
$ws.Range("T3").Value = 0.537552751174421

# Row 4
$ws.Range("E4").Value = 1
$ws.Range("F4").Value = 0.3333333333333333
$ws.Range("G4").Value = 0.05881766666666666
$ws.Range("H4").Value = 0.176453
$ws.Range("M4").Value = 6.58215
$ws.Range("O4").Value = 0.07149673390187571
$ws.Range("P4").Value = 0.07149673390187571
$ws.Range("Q4").Value = 0.38714670465
$ws.Range("R4").Value = 3.48432034185
$ws.Range("S4").Value = 0.07149673390187571
$ws.Range("T4").Value = 0.07149673390187571
